# Fruta / hortaliza, semanal
# The commit reorders the weekly price records (rows 2-19, columns A:R) on
# the sheet into a different date order. Every target row's full content
# (A:R) is identical to some other original row's full content - i.e. this
# edit is a pure re-ordering (permutation) of the existing data rows.
#
# Strategy: snapshot every source row A:R into memory first (so we never
# read a cell after it has already been overwritten), then write the
# snapshots back out in their new row positions according to the mapping
# below (new row number -> original row number that its data came from).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 19

# newRow -> originalRow (row 7 is unchanged and omitted)
$rowMap = @{
    2  = 17
    3  = 9
    4  = 10
    5  = 11
    6  = 15
    8  = 19
    9  = 18
    10 = 3
    11 = 12
    12 = 13
    13 = 8
    14 = 5
    15 = 6
    16 = 4
    17 = 14
    18 = 2
    19 = 16
}

# 1. Snapshot all existing rows (A:R) before making any changes.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = $ws.Range("A" + $r + ":R" + $r).Value()
}

# 2. Write each row's data back into its new location.
foreach ($newRow in $rowMap.Keys) {
    $srcRow = $rowMap[$newRow]
    $ws.Range("A" + $newRow + ":R" + $newRow).Value = $snapshot[$srcRow]
}
